# Apply the "Profile completed except uploading image" edit:
#  - Paragraph 1 ("Forgot Password") becomes three runs reading
#    "SharedPrefrence " + "and parcelable " + "for Username"
#    (the content that used to live in paragraphs 1 and 2 is merged
#    into a single bulleted paragraph).
#  - The old paragraph 2 ("SharedPrefrence for Username", with its
#    proofErr spell-check wrappers) is removed outright - its text
#    now lives in the rewritten paragraph 1.
#  - Paragraph 3 ("Best UI only for pixel 3 , why not on other") is removed.
#  - Paragraph 4 ("Gender radio button") is removed.
#  - The final (empty) paragraph loses its list numbering (<w:numPr>)
#    but keeps the ListParagraph style.

$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

# --- Step 1: rewrite paragraph 1 so it carries the merged sentence as
#     three discrete runs, matching the target markup exactly -------------
$firstXml = '<w:p ' + $wns + ' w14:paraId="507971DF" w14:textId="181DA34C" w:rsidR="00484775" w:rsidRDefault="00654032" w:rsidP="00654032">' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">SharedPrefrence </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">and parcelable </w:t></w:r>' +
  '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>for Username</w:t></w:r>' +
  '</w:p>'
[void]$d.Paragraphs.Item(1).Range.InsertXML($firstXml)

# --- Step 2: delete the three paragraphs that followed (old paragraph 2,
#     which duplicated the sentence above, plus the two unrelated bullets).
#     After step 1 they are still paragraphs 2, 3 and 4; each Delete()
#     removes the whole paragraph (incl. its mark) and shifts the rest up,
#     so repeatedly deleting paragraph 2 removes all three. -----------------
$d.Paragraphs.Item(2).Range.Delete()
$d.Paragraphs.Item(2).Range.Delete()
$d.Paragraphs.Item(2).Range.Delete()

# --- Step 3: strip the numbering from the trailing empty paragraph -------
$lastIndex = $d.Paragraphs.Count
$lastXml = '<w:p ' + $wns + ' w14:paraId="01EACBAE" w14:textId="77777777" w:rsidR="0031618E" w:rsidRPr="00654032" w:rsidRDefault="0031618E" w:rsidP="00654032">' +
  '<w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
  '</w:p>'
[void]$d.Paragraphs.Item($lastIndex).Range.InsertXML($lastXml)
